$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cell H1 "Save", copying the style from G1 (bold/border header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the "Save" values for rows 2-7
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
